$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in columns D (Price) and E (Volume(1h)) are numeric-looking text values
# stored as literal strings in the original workbook (t="inlineStr"). Force the
# cell format to Text before assignment so COM does not auto-convert them into
# real numbers/percentages, then reset the style so no stray number format sticks.
$textCells = @(
    "D2",
    "E2",
    "D3",
    "E3",
    "D4",
    "E4",
    "D5",
    "E5",
    "D6",
    "E6",
    "D7",
    "E7",
    "D8",
    "E8",
    "D9",
    "E9",
    "D10",
    "E10",
    "D11",
    "E11",
    "D12",
    "E12",
    "D13",
    "E13",
    "D14",
    "E14",
    "D15",
    "E15",
    "D16",
    "E16",
    "D17",
    "E17",
    "D18",
    "E18",
    "D19",
    "E19",
    "D20",
    "E20",
    "D21",
    "E21",
    "D22",
    "E22",
    "D23",
    "E23",
    "E24",
    "D25",
    "E25",
    "D26",
    "E26",
    "D27",
    "E27",
    "D39",
    "E39",
    "D40",
    "E40",
    "D41",
    "E41",
    "D42",
    "E42",
    "D43",
    "E43",
    "D44",
    "E44",
    "D45",
    "E45",
    "D46",
    "E46",
    "E47",
    "D48",
    "E48",
    "D49",
    "E49",
    "D50",
    "E50",
    "D51",
    "E51",
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values (coin name / link / price / volume columns).
$ws.Range("D2").Value = "305.26"
$ws.Range("E2").Value = "1.64%"
$ws.Range("D3").Value = "36.34"
$ws.Range("E3").Value = "-4.73%"
$ws.Range("D4").Value = "5.030"
$ws.Range("E4").Value = "0.99%"
$ws.Range("D5").Value = "0.07817"
$ws.Range("E5").Value = "1.29%"
$ws.Range("D6").Value = "2.118"
$ws.Range("E6").Value = "-3.71%"
$ws.Range("D7").Value = "7.907"
$ws.Range("E7").Value = "-0.68%"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "4.096"
$ws.Range("E8").Value = "2.55%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "0.9191"
$ws.Range("E9").Value = "0.35%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "0.09643"
$ws.Range("E10").Value = "5.92%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "0.1866"
$ws.Range("E11").Value = "4.13%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.08645"
$ws.Range("E12").Value = "2.36%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.03510"
$ws.Range("E13").Value = "-0.66%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.09928"
$ws.Range("E14").Value = "-0.07%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "0.001432"
$ws.Range("E15").Value = "-3.49%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "0.005662"
$ws.Range("E16").Value = "0.03%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "3.462"
$ws.Range("E17").Value = "-0.40%"
$ws.Range("D18").Value = "2.522"
$ws.Range("E18").Value = "13.41%"
$ws.Range("D19").Value = "0.3417"
$ws.Range("E19").Value = "-1.32%"
$ws.Range("D20").Value = "0.1345"
$ws.Range("E20").Value = "2.12%"
$ws.Range("D21").Value = "4.763"
$ws.Range("E21").Value = "4.30%"
$ws.Range("D22").Value = "0.2296"
$ws.Range("E22").Value = "2.71%"
$ws.Range("D23").Value = "0.04623"
$ws.Range("E23").Value = "-0.81%"
$ws.Range("E24").Value = "15.08%"
$ws.Range("D25").Value = "0.001233"
$ws.Range("E25").Value = "0.26%"
$ws.Range("D26").Value = "0.0001404"
$ws.Range("E26").Value = "7.88%"
$ws.Range("D27").Value = "0.0004762"
$ws.Range("E27").Value = "0.15%"
$ws.Range("D39").Value = "0.01825"
$ws.Range("E39").Value = "5.06%"
$ws.Range("D40").Value = "0.04748"
$ws.Range("E40").Value = "1.30%"
$ws.Range("D41").Value = "0.007501"
$ws.Range("E41").Value = "-5.00%"
$ws.Range("D42").Value = "0.1401"
$ws.Range("E42").Value = "0.99%"
$ws.Range("D43").Value = "0.007763"
$ws.Range("E43").Value = "1.06%"
$ws.Range("D44").Value = "0.002236"
$ws.Range("E44").Value = "-2.88%"
$ws.Range("D45").Value = "0.01027"
$ws.Range("E45").Value = "5.06%"
$ws.Range("D46").Value = "0.00006270"
$ws.Range("E46").Value = "4.16%"
$ws.Range("E47").Value = "0.16%"
$ws.Range("D48").Value = "0.0005803"
$ws.Range("E48").Value = "0.04%"
$ws.Range("D49").Value = "24.55"
$ws.Range("E49").Value = "181.36%"
$ws.Range("D50").Value = "0.002005"
$ws.Range("E50").Value = "-25.70%"
$ws.Range("D51").Value = "0.00002105"
$ws.Range("E51").Value = "0.16%"

# Restore default (Normal) style on the text cells so no percent/number format
# lingers on them -- matches the unstyled inline-string cells in the source file.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
